$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.464.57"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.616.37"
$ws.Range("E3").Value = "  +5.46%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'238.33"
$ws.Range("E5").Value = "  +1.88%  "
$ws.Range("D6").Value = "'653.26"
$ws.Range("E6").Value = "  +5.09%  "
$ws.Range("E7").Value = "  +3.01%  "
$ws.Range("D8").Value = "'0.404"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "'0.998"
$ws.Range("E10").Value = "  +1.84%  "
$ws.Range("D11").Value = "3.616.27"
$ws.Range("E11").Value = "  +5.42%  "
$ws.Range("D12").Value = "'42.73"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'0.199"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").Value = "'6.33"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "4.309.13"
$ws.Range("E15").Value = "  +6.12%  "
$ws.Range("D16").Value = "95.410.40"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("D17").Value = "'0.0000254"
$ws.Range("E17").Value = "  +2.39%  "
$ws.Range("D18").Value = "3.618.15"
$ws.Range("E18").Value = "  +6.11%  "
$ws.Range("D19").Value = "'7.89"
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").Value = "'12.87"
$ws.Range("E20").Value = "  +9.81%  "
$ws.Range("D21").Value = "'17.97"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "'3.63"
$ws.Range("E22").Value = "  +6.44%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").Value = "'508.31"
$ws.Range("E24").Value = "  +1.29%  "
$ws.Range("D25").Value = "'0.0000198"
$ws.Range("E25").Value = "  +6.60%  "
$ws.Range("D26").Value = "'6.59"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "'96.18"
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("D28").Value = "'12.66"
$ws.Range("E28").Value = "  +4.33%  "
$ws.Range("D29").Value = "3.795.06"
$ws.Range("E29").Value = "  +5.11%  "
$ws.Range("D30").Value = "'3.15"
$ws.Range("E30").Value = "  +13.88%  "
$ws.Range("D31").Value = "'11.35"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("E34").Value = "  +1.07%  "
$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("D36").Value = "'32.17"
$ws.Range("E36").Value = "  +9.41%  "
$ws.Range("D37").Value = "'0.561"
$ws.Range("E37").Value = "  +0.69%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D38").Value = "'8.17"
$ws.Range("E38").Value = "  +8.60%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D39").Value = "'571.84"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'1.47"
$ws.Range("E40").Value = "  +3.71%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'0.930"
$ws.Range("E42").Value = "  +2.36%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.150"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "'35.05"
$ws.Range("E44").Value = "  +41.30%  "
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").Value = "'1.72"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'23.74"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").Value = "'5.68"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "'2.24"
$ws.Range("E48").Value = "  +4.64%  "
$ws.Range("D49").Value = "'0.0414"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("E50").Value = "  -5.39%  "
$ws.Range("D51").Value = "'53.78"
$ws.Range("E51").Value = "  +0.59%  "
